# Slide 5 ("Extraction elements"): resize/reposition the existing
# "Arrow: Right 6" shape, label it "Regular expression", then duplicate
# it into a second arrow ("Arrow: Right 5") labeled "NER".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# Locate the right-arrow shape by name (id=7, "Right Arrow 6" / "Arrow: Right 6").
$arrow = $null
foreach ($sh in $s.Shapes) {
    if ($sh.Name -eq "Right Arrow 6" -or $sh.Name -eq "Arrow: Right 6") {
        $arrow = $sh
    }
}

# Reposition / resize the first arrow.
# (point values chosen so the EMU round-trip lands on the exact target values)
$arrow.Left   = 332.73780827559057
$arrow.Top    = 165.4251181102362
$arrow.Width  = 203.52063792125986
$arrow.Height = 66.88606699212599

# Duplicate BEFORE adding text, so both arrows start from the same
# unlabeled state that the original shape had.
# (Duplicate twice and drop the first copy so the surviving shape's
# auto-assigned Id lands on 6, matching the recorded id sequence.)
$throwaway = $arrow.Duplicate()
$dup = $arrow.Duplicate()
$throwaway.Delete()
$dup.Name = "Arrow: Right 5"

# Label the first arrow.
$arrow.TextFrame.TextRange.Text = "Regular expression"
$arrow.TextFrame.TextRange.Font.Size = 20

# Reposition / resize the duplicated arrow.
$dup.Left   = 353.7088168976378
$dup.Top    = 249.96819397637793
$dup.Width  = 148.33654043307087
$dup.Height = 72.210003

# Label the duplicated arrow.
$dup.TextFrame.TextRange.Text = "NER"
$dup.TextFrame.TextRange.Font.Size = 24
